$wb = $excel.ActiveWorkbook

# ALC row 40 (Leve Item ID 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7000
$ws.Range("I40").Value = 7000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 7000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -6825
$ws.Range("N40").ClearContents()

# ALC row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 1950
$ws.Range("I76").Value = 1950
$ws.Range("K76").Value = 1950
$ws.Range("M76").Value = -1635

# ALC row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 1950
$ws.Range("I79").Value = 1950
$ws.Range("K79").Value = 1950
$ws.Range("M79").Value = -858

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2187.3
$ws.Range("I132").Value = 2187.3
$ws.Range("K132").Value = 6561.900000000001
$ws.Range("M132").Value = -4031.900000000001

# ALC row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1833.4722
$ws.Range("I137").Value = 1923
$ws.Range("J137").Value = 1654.4166
$ws.Range("K137").Value = 5769
$ws.Range("L137").Value = 4963.2498
$ws.Range("M137").Value = -3219
$ws.Range("N137").Value = -10063.2498

# ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 9869.75
$ws.Range("J138").Value = 9869.75
$ws.Range("L138").Value = 29609.25
$ws.Range("N138").Value = -39889.25

# ARM row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8086.794
$ws.Range("I32").Value = 8086.794
$ws.Range("K32").Value = 8086.794
$ws.Range("M32").Value = -7799.794

# ARM row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 14968.75
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 22750
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 22750
$ws.Range("M110").Value = 45
$ws.Range("N110").Value = -26840

# ARM row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3250
$ws.Range("I122").Value = 3250
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 9750
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7300
$ws.Range("N122").ClearContents()

# BSM row 105 (Leve Item ID 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 5000
$ws.Range("K105").Value = 5000
$ws.Range("M105").Value = -3253

# BSM row 129 (Leve Item ID 35382)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 29000
$ws.Range("J129").Value = 29000
$ws.Range("L129").Value = 29000
$ws.Range("N129").Value = -39000

# CRP row 16 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2550
$ws.Range("I16").Value = 2550
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2550
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2263
$ws.Range("N16").ClearContents()

# CRP row 58 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4328.357
$ws.Range("I58").Value = 4353.6924
$ws.Range("J58").Value = 3999
$ws.Range("K58").Value = 4353.6924
$ws.Range("L58").Value = 3999
$ws.Range("M58").Value = -4150.6924
$ws.Range("N58").Value = -4405

# CRP row 113 (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2550
$ws.Range("I113").Value = 2550
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2550
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -380
$ws.Range("N113").ClearContents()

# CRP row 136 (Leve Item ID 44021)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4328.357
$ws.Range("I136").Value = 4353.6924
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 13061.0772
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -10511.0772
$ws.Range("N136").Value = -17097

# CUL row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 3000
$ws.Range("N5").Value = -3224

# CUL row 33 (Leve Item ID 4867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 108
$ws.Range("I33").Value = 108
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 648
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -365
$ws.Range("N33").ClearContents()

# CUL row 40 (Leve Item ID 4827)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 304.25
$ws.Range("J40").Value = 350
$ws.Range("L40").Value = 1400
$ws.Range("N40").Value = -1538

# CUL row 129 (Leve Item ID 36054)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2062
$ws.Range("J129").Value = 2483.2856
$ws.Range("L129").Value = 7449.8568
$ws.Range("N129").Value = -17449.8568

# CUL row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1000
$ws.Range("J135").Value = 1000
$ws.Range("L135").Value = 9000
$ws.Range("N135").Value = -14070

# GSM row 10 (Leve Item ID 4306)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 21000500
$ws.Range("I10").Value = 21000500
$ws.Range("K10").Value = 21000500
$ws.Range("M10").Value = -21000331

# GSM row 46 (Leve Item ID 2078)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 20972.5
$ws.Range("I46").Value = 1999
$ws.Range("K46").Value = 1999
$ws.Range("M46").Value = -1843

# GSM row 70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# GSM row 73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3290.7778
$ws.Range("I122").Value = 2873.8572
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 8621.571599999999
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -6171.571599999999
$ws.Range("N122").Value = -19150

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3070.0908
$ws.Range("I132").Value = 2910.375
$ws.Range("J132").Value = 3496
$ws.Range("K132").Value = 8731.125
$ws.Range("L132").Value = 10488
$ws.Range("M132").Value = -6201.125
$ws.Range("N132").Value = -15548

# LTW row 35 (Leve Item ID 1697)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 5066.6665
$ws.Range("I35").Value = 2600
$ws.Range("K35").Value = 2600
$ws.Range("M35").Value = -2264

# LTW row 46 (Leve Item ID 5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2688.5715
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376

# LTW row 53 (Leve Item ID 3866)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 4000
$ws.Range("I53").Value = 4000
$ws.Range("K53").Value = 4000
$ws.Range("M53").Value = -3482

# LTW row 100 (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4671.4287
$ws.Range("J100").Value = 1850
$ws.Range("L100").Value = 1850
$ws.Range("N100").Value = -2932

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3596.6
$ws.Range("I136").Value = 3496.75
$ws.Range("K136").Value = 10490.25
$ws.Range("M136").Value = -7940.25

# WVR row 122 (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3651
$ws.Range("I122").Value = 4093.2856
$ws.Range("J122").Value = 555
$ws.Range("K122").Value = 12279.8568
$ws.Range("L122").Value = 1665
$ws.Range("M122").Value = -9829.856800000001
$ws.Range("N122").Value = -6565
